# Weekly update: insert a new data row at row 69 (pushing existing rows
# 69-113 down to 70-114) for the "Hortaliza, Vega Monumental Concepción -
# Alcachofa" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 69; this shifts rows
# 69..113 down to 70..114 and extends the sheet dimension automatically.
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the new weekly record.
$ws.Cells.Item(69, 1).Value2 = 11
$ws.Cells.Item(69, 2).Value2 = 'Vega Monumental Concepción'
$ws.Cells.Item(69, 3).Value2 = 'Bíobío'
$ws.Cells.Item(69, 4).Value2 = 45161
$ws.Cells.Item(69, 5).Value2 = 8
$ws.Cells.Item(69, 6).Value2 = 100112013
$ws.Cells.Item(69, 7).Value2 = 'Alcachofa'
$ws.Cells.Item(69, 8).Value2 = 'Argentina(o)'
$ws.Cells.Item(69, 9).Value2 = 'Primera'
$ws.Cells.Item(69, 10).Value2 = 50
$ws.Cells.Item(69, 11).Value2 = 13000
$ws.Cells.Item(69, 12).Value2 = 13000
$ws.Cells.Item(69, 13).Value2 = 13000
$ws.Cells.Item(69, 14).Value2 = '$/caja 50 unidades'
$ws.Cells.Item(69, 15).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(69, 16).Value2 = 260
$ws.Cells.Item(69, 17).Value2 = 50
$ws.Cells.Item(69, 18).Value2 = 'Hortaliza'
